$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 61; this shifts former rows 61-88 down to 62-89
# and extends the used range to A1:T89.
$ws.Rows(61).Insert()

# Fill in the new row 61 with the new market entry.
$ws.Cells.Item(61, 1).Value = 4
$ws.Cells.Item(61, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(61, 3).Value = "Los Lagos"
$ws.Cells.Item(61, 4).Value = 44719
$ws.Cells.Item(61, 5).Value = 10
$ws.Cells.Item(61, 6).Value = "Fruta"
$ws.Cells.Item(61, 7).Value = 100104
$ws.Cells.Item(61, 8).Value = "Frutos de pepita"
$ws.Cells.Item(61, 9).Value = 100104003
$ws.Cells.Item(61, 10).Value = "Membrillo"
$ws.Cells.Item(61, 11).Value = "Champion"
$ws.Cells.Item(61, 12).Value = "Primera"
$ws.Cells.Item(61, 13).Value = 400
$ws.Cells.Item(61, 14).Value = 13000
$ws.Cells.Item(61, 15).Value = 14000
$ws.Cells.Item(61, 16).Value = 13500
$ws.Cells.Item(61, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(61, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(61, 19).Value = 750
$ws.Cells.Item(61, 20).Value = 18
